$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "State" (B) and "District"/other (C) columns of data, entered in the
# same order the original author typed them (column B top-to-bottom, then
# column C top-to-bottom but C3 before C2) so the shared-string table comes
# out in the same order as the authored workbook.
$ws.Range("B1").Value = "Patna"
$ws.Range("B2").Value = "Bihar"
$ws.Range("B3").Value = "India"
$ws.Range("C1").Value = "Bihar"
$ws.Range("C3").Value = "Amnour"
$ws.Range("C2").Value = "Saran"

$ws.Range("C1").Select() | Out-Null
